$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D/E column cells to match the refreshed crypto snapshot.
# D-column price strings must stay as TEXT (the sheet stores prices as
# strings, some with "." thousands separators) - temporarily force the
# cell to Text format so Excel does not reinterpret the digits as a
# number, then restore the original (default) style so no stray
# number-format/quote-prefix style is left behind.

$priceCell = $ws.Range("D2")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '63.741.64'
$priceCell.Style = $origStyle
$ws.Range("E2").Value = '  +2.83%  '

$priceCell = $ws.Range("D3")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '3.130.23'
$priceCell.Style = $origStyle
$ws.Range("E3").Value = '  +1.48%  '

$ws.Range("E4").Value = '  -0.19%  '

$priceCell = $ws.Range("D5")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '589.97'
$priceCell.Style = $origStyle
$ws.Range("E5").Value = '  +1.74%  '

$priceCell = $ws.Range("D6")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '146.72'
$priceCell.Style = $origStyle
$ws.Range("E6").Value = '  +2.81%  '

$ws.Range("E7").Value = '  -0.05%  '

$priceCell = $ws.Range("D8")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '3.124.16'
$priceCell.Style = $origStyle
$ws.Range("E8").Value = '  +1.74%  '

$ws.Range("E9").Value = '  +1.63%  '

$ws.Range("E10").Value = '  +16.21%  '

$priceCell = $ws.Range("D11")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '5.74'
$priceCell.Style = $origStyle
$ws.Range("E11").Value = '  +3.96%  '

$priceCell = $ws.Range("D12")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.469'
$priceCell.Style = $origStyle
$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("E13").Value = '  +5.55%  '

$priceCell = $ws.Range("D14")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '36.40'
$priceCell.Style = $origStyle
$ws.Range("E14").Value = '  +3.26%  '

$ws.Range("E15").Value = '  -0.54%  '

$priceCell = $ws.Range("D16")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '3.648.72'
$priceCell.Style = $origStyle
$ws.Range("E16").Value = '  +1.49%  '

$priceCell = $ws.Range("D17")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '7.19'
$priceCell.Style = $origStyle
$ws.Range("E17").Value = '  -0.95%  '

$priceCell = $ws.Range("D18")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '63.648.27'
$priceCell.Style = $origStyle
$ws.Range("E18").Value = '  +2.78%  '

$priceCell = $ws.Range("D19")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '3.124.76'
$priceCell.Style = $origStyle
$ws.Range("E19").Value = '  +1.26%  '

$priceCell = $ws.Range("D20")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '465.70'
$priceCell.Style = $origStyle
$ws.Range("E20").Value = '  +3.88%  '

$priceCell = $ws.Range("D21")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '14.44'
$priceCell.Style = $origStyle
$ws.Range("E21").Value = '  +3.73%  '

$priceCell = $ws.Range("D22")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.737'
$priceCell.Style = $origStyle
$ws.Range("E22").Value = '  +0.92%  '

$priceCell = $ws.Range("D23")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '7.57'
$priceCell.Style = $origStyle
$ws.Range("E23").Value = '  +1.41%  '

$priceCell = $ws.Range("D24")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '13.32'
$priceCell.Style = $origStyle
$ws.Range("E24").Value = '  -3.61%  '

$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("E26").Value = '  -0.01%  '

$priceCell = $ws.Range("D27")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '9.02'
$priceCell.Style = $origStyle
$ws.Range("E27").Value = '  +11.39%  '

$ws.Range("E28").Value = '  +2.20%  '

$ws.Range("E29").Value = '  -1.05%  '

$ws.Range("E30").Value = '  -0.15%  '

$priceCell = $ws.Range("D31")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '6.92'
$priceCell.Style = $origStyle
$ws.Range("E31").Value = '  +3.13%  '

$priceCell = $ws.Range("D32")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '27.28'
$priceCell.Style = $origStyle
$ws.Range("E32").Value = '  +2.12%  '

$priceCell = $ws.Range("D33")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.110'
$priceCell.Style = $origStyle
$ws.Range("E33").Value = '  -2.84%  '

$priceCell = $ws.Range("D34")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.0₃0865'
$priceCell.Style = $origStyle
$ws.Range("E34").Value = '  +8.12%  '

$priceCell = $ws.Range("D35")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '2.38'
$priceCell.Style = $origStyle
$ws.Range("E35").Value = '  +8.92%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("E37").Value = '  +12.83%  '

$priceCell = $ws.Range("D38")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '6.12'
$priceCell.Style = $origStyle
$ws.Range("E38").Value = '  +1.37%  '

$priceCell = $ws.Range("D39")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '51.03'
$priceCell.Style = $origStyle
$ws.Range("E39").Value = '  +1.24%  '

$priceCell = $ws.Range("D40")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '450.45'
$priceCell.Style = $origStyle
$ws.Range("E40").Value = '  +4.59%  '

$priceCell = $ws.Range("D41")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '8.80'
$priceCell.Style = $origStyle
$ws.Range("E41").Value = '  -0.36%  '

$priceCell = $ws.Range("D42")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.0374'
$priceCell.Style = $origStyle
$ws.Range("E42").Value = '  +0.77%  '

$priceCell = $ws.Range("D43")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '2.899.25'
$priceCell.Style = $origStyle
$ws.Range("E43").Value = '  +3.79%  '

$priceCell = $ws.Range("D44")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '0.278'
$priceCell.Style = $origStyle
$ws.Range("E44").Value = '  +3.52%  '

$ws.Range("E45").Value = '  +2.34%  '

$priceCell = $ws.Range("D46")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '2.22'
$priceCell.Style = $origStyle
$ws.Range("E46").Value = '  +5.95%  '

$priceCell = $ws.Range("D47")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '35.82'
$priceCell.Style = $origStyle
$ws.Range("E47").Value = '  +1.10%  '

$priceCell = $ws.Range("D48")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '126.52'
$priceCell.Style = $origStyle
$ws.Range("E48").Value = '  +1.22%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E50").Value = '  +0.42%  '

$priceCell = $ws.Range("D51")
$origStyle = $priceCell.Style
$priceCell.NumberFormat = "@"
$priceCell.Value = '24.84'
$priceCell.Style = $origStyle
$ws.Range("E51").Value = '  +3.39%  '
